$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The roster rows for Cam Thomas (row 5) and Seth Curry (row 6) were
# reordered (Seth Curry now listed before Cam Thomas). Swap every
# per-player field between the two rows to reflect the new ordering.

# No.
$ws.Range("B5").Value = 30
$ws.Range("B6").Value = 24

# Player
$ws.Range("C5").Value = "Seth Curry"
$ws.Range("C6").Value = "Cam Thomas"

# Ht
$ws.Range("E5").Value = "6-2"
$ws.Range("E6").Value = "6-4"

# Wt
$ws.Range("F5").Value = 185
$ws.Range("F6").Value = 210

# Birth Date
$ws.Range("G5").Value = "August 23, 1990"
$ws.Range("G6").Value = "October 13, 2001"

# Unnamed: 6 (country code)
$ws.Range("H5").Value = "us"
$ws.Range("H6").Value = "jp"

# Exp (format as text first so "8"/"1" stay text like the rest of the column,
# then restore the cell style so no stray number-format sticks around)
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = "8"
$ws.Range("I5").Style = "Normal"

$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = "1"
$ws.Range("I6").Style = "Normal"

# College
$ws.Range("J5").Value = "Liberty, Duke"
$ws.Range("J6").Value = "LSU"

# bbref url
$ws.Range("K5").Value = "https://www.basketball-reference.com/players/c/curryse01.html"
$ws.Range("K6").Value = "https://www.basketball-reference.com/players/t/thomaca02.html"
